# Edit slide 1: update title + rewrite/extend the bullet paragraphs in the
# content placeholder, per the target diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title: "Lebron the GOAT of NBA" -> "Test" ---------------------------
# Replace via a full-length Characters() sub-range (rather than assigning
# TextRange.Text directly) so the existing run isn't given a synthesized
# a:rPr/lang attribute it didn't have before.
$titleShape = $s.Shapes.Item(1)
$titleTr = $titleShape.TextFrame.TextRange
($titleTr.Characters(1, $titleTr.Text.Length)).Text = "Test"

# --- Content placeholder: rewrite existing bullets, add a new one -------
$contentShape = $s.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange

# Helper: build a sub-TextRange covering the Nth paragraph (0-based, split
# on the `r paragraph separator PowerPoint uses in TextRange.Text) so we can
# replace its text in place without disturbing a:pPr / run formatting on
# the surrounding paragraphs.
function Get-ParaRange($textRange, $paraIndex) {
    $full = $textRange.Text
    $parts = $full -split "`r"
    $start = 1
    for ($i = 0; $i -lt $paraIndex; $i++) {
        $start += $parts[$i].Length + 1
    }
    $len = $parts[$paraIndex].Length
    return $textRange.Characters($start, $len)
}

# Paragraph layout before edit:
#   0: "" (empty paragraph)
#   1: "Lebron is considered for many the greatest basketball player of all time"
#   2: "There people who compare him to Michael Jordan, but everyone knows who is the goat"

(Get-ParaRange $tr 1).Text = "Discussion on LeBron James being a top basketball player who has played for multiple teams and won titles during the pandemic bubble."
(Get-ParaRange $tr 2).Text = "Mention of his career including stints with Cleveland Cavaliers, Miami Heat, and current team LA Lakers."

# Append a brand-new 4th paragraph (inherits the preceding paragraph's
# pPr/defRPr formatting, matching the diff's new <a:p><a:pPr><a:defRPr
# sz="1600"/></a:pPr> block). Discard the returned TextRange so it isn't
# echoed to the output stream.
$null = $tr.InsertAfter("`rComparison with Jordan, asserting LeBron's greatness in basketball.")
